$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Split " in Power Engineering" into " in Power Engineerin" + "g", with a
#    _GoBack bookmark inserted between them (this also relocates the existing
#    _GoBack bookmark away from its old spot near "both the utilities...").
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Power Engineerin", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $rng.End
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ---------------------------------------------------------------------------
# 2) Renumber the sole document comment from id 0 to id 1. The Word object
#    model does not expose the raw comment id, but adding a brand new comment
#    anchored to the exact same scope (and then deleting the original one)
#    causes the surviving comment to be written out as id "1" (ids are
#    allocated in order and the old id 0 is freed up by the delete).
# ---------------------------------------------------------------------------
$origComment = $d.Comments(1)
$commentScope = $origComment.Scope
$lb = [char]11
$commentText = 'This paragraph needs to be reorganized/simplified but I' + [char]8217 + 'm not exactly sure how to make it better.' + $lb + `
    'Joined DLC to improve and finish project' + $lb + `
    'Describe why the project was useful' + $lb + `
    'Describe the progress I' + [char]8217 + 've made on improving the project' + $lb + `
    'Describe why the project was unique' + $lb + `
    'Published a paper on the project.'
$newComment = $d.Comments.Add($commentScope, $commentText)
$d.Comments(1).Delete() | Out-Null
$renumbered = $d.Comments(1)
$renumbered.Author = "Morgenstern, Carl W"
$renumbered.Initial = "MCW"

# ---------------------------------------------------------------------------
# 3) Rewrite the tutoring paragraph.
# ---------------------------------------------------------------------------
$old1 = "During our tutoring time I am able to impart some of my young wisdom and life advice. By building a bond with him I hope to nudge"
$new1 = "Before our tutoring sessions and we make dinner together, and during this time we talk about things outside of school. I try to bring up topics so that he might begin to think of issues that are bigger than his next pair of shoes, such as gender roles, kneeling for the national anthem, and life after high school. I hope the bond I am trying to build with him nudges"
$r1 = $d.Content
$r1.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

$old2 = "in the direction of enjoying school, building empathy, and fulfilling his dreams of going to a prestigious east coast college. "
$new2 = "in the right direction, and he can evolve from a surly teenager to a thoughtful and involved human being. "
$r2 = $d.Content
$r2.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Append a closing sentence after the final paragraph's last sentence.
# ---------------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("epitomize all the characteristics of a Graduate Research Fellow.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$endPos = $r3.End
$tailRng = $d.Range($endPos, $endPos)
$tailRng.InsertAfter(" Thank you for your consideration.")
